$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.123.11'
$ws.Range('E2').Value = '  +0.94%  '
$ws.Range('D3').Value = '3.258.83'
$ws.Range('E3').Value = '  +0.12%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '581.58'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '184.30'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.64%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.595'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.68%  '
$ws.Range('E9').Value = '  +0.79%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.64'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.35%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.419'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.92%  '
$ws.Range('D12').Value = '3.809.11'
$ws.Range('E12').Value = '  -0.21%  '
$ws.Range('E13').Value = '  -0.41%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.25'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.49%  '
$ws.Range('D15').Value = '68.136.51'
$ws.Range('E15').Value = '  +1.02%  '
$ws.Range('D17').Value = '3.227.52'
$ws.Range('E17').Value = '  -1.42%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.82'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.32%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.53'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '391.33'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.48%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.73'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '71.46'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.65%  '
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.520'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.69%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000120'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.22%  '
$ws.Range('E26').Value = '  +3.75%  '
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.997'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.85%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.98'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.10%  '
$ws.Range('E30').Value = '  +0.42%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '23.02'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.34%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.18'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.86%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.29'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.96%  '
$ws.Range('B34').Value = 'USDe'
$ws.Range('C34').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.998'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '164.90'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.50'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.71%  '
$ws.Range('E37').Value = '  +3.96%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.822'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.84%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '26.70'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.03%  '
$ws.Range('E40').Value = '  -0.55%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.56'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.08%  '
$ws.Range('E42').Value = '  -3.14%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '41.35'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.41%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0683'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.85%  '
$ws.Range('D45').Value = '2.629.85'
$ws.Range('E45').Value = '  -2.45%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '341.17'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.81%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '24.82'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.39%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0281'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.34%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.33'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.09%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '31.64'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.19%  '
$ws.Range('E51').Value = '  -0.21%  '
